$p = $ppt.ActivePresentation

# --- Change 1: slide 3 ------------------------------------------------
# Merge the stray " " run with the following "is that it provides..."
# run into a single run (same formatting on both, so re-writing the
# combined character range collapses them into one <a:r>).
$slide3 = $p.Slides.Item(3)
$titleShape = $slide3.Shapes.Item(1)
$combined = $titleShape.TextFrame.TextRange.Characters(90, 63)
$combined.Text = " is that it provides a uniform interface to start unites. This "

# --- Change 2: slide 6 --------------------------------------------------
# Give the ctrTitle placeholder an explicit position/size (it previously
# had an empty <p:spPr/>, inheriting from the layout).
$slide6 = $p.Slides.Item(6)
$ctrTitleShape = $slide6.Shapes.Item(1)
$ctrTitleShape.Left = 42
$ctrTitleShape.Top = 168
$ctrTitleShape.Width = 612
$ctrTitleShape.Height = 115.75
